$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing columns A-D shift to B-E
$ws.Columns("A").Insert()

# Updated Cypher query text. Set the stat-query (C2) before the trials-filter
# query (B2) so new shared strings are appended in the same order as the
# target workbook (StatQuery text, then trials-filter text, then the new
# TabName/CasesTab labels).
$si6 = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@

$si7 = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "AMERICAN_INDIAN_OR_ALASKA_NATIVE"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@

$ws.Range("C2").Value = $si6
$ws.Range("B2").Value = $si7

# New "TabName" column A header/value
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# New column A is narrow (best-fit style); the rest keep their original widths
# (column A is new, so it has no pre-existing width to inherit).
$ws.Columns("A").ColumnWidth = 8.3

# Row 2 needs to grow to fit the longer query text
$ws.Rows("2").RowHeight = 174

# Selection moves to B4 in the saved file
$ws.Range("B4").Select() | Out-Null

